$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Fix the gender recode: `1` = 0, `2` = 1  ->  `1` = 1, `2` = 0
#    (done with precise single-character Range edits so the existing
#    per-token syntax-highlighting runs -- DecValTok, StringTok, etc. --
#    stay intact, matching the target diff which only flips the two digits.)
# ---------------------------------------------------------------------------
$recodePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.IndexOf("gender = recode(gender, ``1``") -ge 0) {
        $recodePara = $cand
        break
    }
}

$pText = $recodePara.Range.Text
$pStart = $recodePara.Range.Start

$needleA = "gender = recode(gender, ``1`` = "
$idxA = $pText.IndexOf($needleA)
$zeroPos = $idxA + $needleA.Length
$rZero = $d.Range($pStart + $zeroPos, $pStart + $zeroPos + 1)
$rZero.Text = "1"

$pText2 = $recodePara.Range.Text
$needleB = "``2`` = "
$idxB = $pText2.IndexOf($needleB, $zeroPos)
$onePos = $idxB + $needleB.Length
$rOne = $d.Range($pStart + $onePos, $pStart + $onePos + 1)
$rOne.Text = "0"

# ---------------------------------------------------------------------------
# 2. Swap the gender column (0 <-> 1) in the printed tibble preview
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "##  1 CAN           492.  0.93       NA         NA      0          NA            3",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##  1 CAN           492.  0.93       NA         NA      1          NA            3", 2) | Out-Null

$d.Content.Find.Execute(
    "##  2 CAN           394. -0.78        0         NA      0          NA            3",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##  2 CAN           394. -0.78        0         NA      1          NA            3", 2) | Out-Null

$d.Content.Find.Execute(
    "##  3 CAN           390. -1.3         0          1      1           2            2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##  3 CAN           390. -1.3         0          1      0           2            2", 2) | Out-Null

$d.Content.Find.Execute(
    "##  4 CAN           504.  0.56        0          2      0           2            3",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##  4 CAN           504.  0.56        0          2      1           2            3", 2) | Out-Null

$d.Content.Find.Execute(
    "##  5 CAN           466. -0.03        0          3      1           1           NA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##  5 CAN           466. -0.03        0          3      0           1           NA", 2) | Out-Null

$d.Content.Find.Execute(
    "##  6 CAN           398.  0.74        0          1      0           2            2",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##  6 CAN           398.  0.74        0          1      1           2            2", 2) | Out-Null

$d.Content.Find.Execute(
    "##  7 CAN           404. NA          NA         NA      0          NA           NA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##  7 CAN           404. NA          NA         NA      1          NA           NA", 2) | Out-Null

$d.Content.Find.Execute(
    "##  8 CAN           406. -2.58        0          4      0           2           NA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##  8 CAN           406. -2.58        0          4      1           2           NA", 2) | Out-Null

$d.Content.Find.Execute(
    "##  9 CAN           609.  0.88        0          4      1           1           NA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##  9 CAN           609.  0.88        0          4      0           1           NA", 2) | Out-Null

$d.Content.Find.Execute(
    "## 10 CAN           452.  0.44        0          1      0           2           NA",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## 10 CAN           452.  0.44        0          1      1           2           NA", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Rewrite the "zero-order correlation" sentence
# ---------------------------------------------------------------------------
$oldSentence = "of a zero-order correlation, which we can think of as a " + [char]0x201C + "normal" + [char]0x201D + " correlation (i.e., Pearson" + [char]0x2019 + "s Product Moment Correlation)."
$newSentence = "of the multiple correlation coefficient; when the part correlation is squared, it carves out the variance explained (i.e., the multiple R^2) by the unique predictor."
$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Insert two new Body Text paragraphs after the rewritten sentence
#    and before the "# Select the variables in the model" source block.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("A part correlation")) {
        $targetPara = $cand
        break
    }
}

$targetPara.Range.InsertParagraphAfter() | Out-Null

$newPara1 = $targetPara.Next()
$newPara1.Range.Text = "Partial is like a relative frequency. It is out of 100% of multiple R^2."
$newPara1.Style = "Body Text"

$newPara1.Range.InsertParagraphAfter() | Out-Null
$newPara2 = $newPara1.Next()
$newPara2.Range.Text = "Whereas the part correlation carves out a part of the multiple R^2. For example, if the R^2 is 0.241 and the part correlation for SES is .483, then the squared value of the part correlation is 0.231, which is in the same units as multiple R^2."
$newPara2.Style = "Body Text"

# ---------------------------------------------------------------------------
# 5. Fix the $estimate correlation matrix (spacing + sign flips on gender col/row)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "##            math_score          ses    language  enjoy_math      gender",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##             math_score          ses     language  enjoy_math       gender", 2) | Out-Null

$d.Content.Find.Execute(
    "## math_score 1.00000000  0.482515108 0.031543578  0.06892420 0.038743018",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## math_score  1.00000000  0.482515108  0.031543578  0.06892420 -0.038743018", 2) | Out-Null

$d.Content.Find.Execute(
    "## ses        0.47969178  1.000000000 0.040210691 -0.12373319 0.007309905",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## ses         0.47969178  1.000000000  0.040210691 -0.12373319 -0.007309905", 2) | Out-Null

$d.Content.Find.Execute(
    "## language   0.03606365  0.046243313 1.000000000  0.03294749 0.006582796",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## language    0.03606365  0.046243313  1.000000000  0.03294749 -0.006582796", 2) | Out-Null

$d.Content.Find.Execute(
    "## enjoy_math 0.07812423 -0.141074622 0.032664616  1.00000000 0.059851654",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## enjoy_math  0.07812423 -0.141074622  0.032664616  1.00000000 -0.059851654", 2) | Out-Null

$d.Content.Find.Execute(
    "## gender     0.04430790  0.008409073 0.006584751  0.06038789 1.000000000",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## gender     -0.04430790 -0.008409073 -0.006584751 -0.06038789  1.000000000", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Fix the $statistic matrix (spacing + sign flips on gender col/row)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "##            math_score        ses language enjoy_math    gender",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "##            math_score        ses  language enjoy_math     gender", 2) | Out-Null

$d.Content.Find.Execute(
    "## math_score   0.000000 112.461929 6.442733  14.104210  7.915214",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## math_score   0.000000 112.461929  6.442733  14.104210  -7.915214", 2) | Out-Null

$d.Content.Find.Execute(
    "## ses        111.606456   0.000000 8.215537 -25.455372  1.492335",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## ses        111.606456   0.000000  8.215537 -25.455372  -1.492335", 2) | Out-Null

$d.Content.Find.Execute(
    "## language     7.367079   9.450544 0.000000   6.729784  1.343887",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## language     7.367079   9.450544  0.000000   6.729784  -1.343887", 2) | Out-Null

$d.Content.Find.Execute(
    "## enjoy_math  15.997723 -29.090901 6.671944   0.000000 12.240480",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## enjoy_math  15.997723 -29.090901  6.671944   0.000000 -12.240480", 2) | Out-Null

$d.Content.Find.Execute(
    "## gender       9.054217   1.716748 1.344286  12.350546  0.000000",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## gender      -9.054217  -1.716748 -1.344286 -12.350546   0.000000", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7. Fix the lm() summary coefficients (Intercept & gender)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "## (Intercept) 449.3502     1.1874 378.429  < 2e-16 ***",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## (Intercept) 456.6267     1.2423 367.579  < 2e-16 ***", 2) | Out-Null

$d.Content.Find.Execute(
    "## gender        7.2765     0.8016   9.077  < 2e-16 ***",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## gender       -7.2765     0.8016  -9.077  < 2e-16 ***", 2) | Out-Null

# ---------------------------------------------------------------------------
# 8. Fix the confint() output (Intercept & gender)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "## (Intercept) 447.022873 451.677573",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## (Intercept) 454.191858 459.061551", 2) | Out-Null

$d.Content.Find.Execute(
    "## gender        5.705252   8.847711",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## gender       -8.847711  -5.705252", 2) | Out-Null

# ---------------------------------------------------------------------------
# 9. Update the closing interpretation paragraph
# ---------------------------------------------------------------------------
$oldClose = "Partial correlations were conducted to examine the unique contributions to the overall variance in math scores. Results indicated that standardized socioeconomic status explained 23.01% of the variance in math achievement, speaking a language different from the test explained 0.13%, student enjoyment of math explained 0.61%, and being female explained 0.2%. Therefore, only socioeconomic status explained a meaningful amount of the variance in math achievement."
$newClose = "Part correlations were conducted to examine the unique contributions to the overall variance in math scores. Part correlations were squared to determine the unique variance. Results indicated that standardized socioeconomic status explained 23.01% of the variance in math achievement, speaking a language different from the test explained 0.13%, student enjoyment of math explained 0.61%, and being female explained 0.2%. Therefore, only socioeconomic status explained a meaningful amount of the variance in math achievement."
$d.Content.Find.Execute($oldClose, $true, $false, $false, $false, $false, $true, 1, $false, $newClose, 2) | Out-Null

Write-Output "done"
